$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F
$ws.Range("F1").Value = "sd_Y"

# New sd_Y values for rows 2-17
$values = @(
    0.034729044126633,
    0.17362021571538,
    0.163491373963925,
    0.0959014601037375,
    0,
    0.181276102183571,
    0.118599994384477,
    0.0687493707269867,
    0,
    0.24833355257902,
    0.2531483122925,
    0.00337751101205847,
    0.0307688497355245,
    0.235441512787457,
    0.120587490821949,
    0
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
